# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted into the data table right
# before the existing row 112, pushing every subsequent row down by one
# (old row 112 -> new row 113, ..., old row 174 -> new row 175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 112; Excel copies formatting
# (including the date number format on column D) from the row above,
# which matches the surrounding rows.
$ws.Rows(112).Insert()

# Populate the freshly inserted row 112 with the new record.
$ws.Cells.Item(112, 1).Value2 = 3
$ws.Cells.Item(112, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(112, 3).Value2 = "Coquimbo"
$ws.Cells.Item(112, 4).Value2 = 44572
$ws.Cells.Item(112, 5).Value2 = 5
$ws.Cells.Item(112, 6).Value2 = "Fruta"
$ws.Cells.Item(112, 7).Value2 = 100101
$ws.Cells.Item(112, 8).Value2 = "Berries"
$ws.Cells.Item(112, 9).Value2 = 100101001
$ws.Cells.Item(112, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(112, 11).Value2 = "Sin especificar"
$ws.Cells.Item(112, 12).Value2 = "Primera"
$ws.Cells.Item(112, 13).Value2 = 30
$ws.Cells.Item(112, 14).Value2 = 5000
$ws.Cells.Item(112, 15).Value2 = 5000
$ws.Cells.Item(112, 16).Value2 = 5000
$ws.Cells.Item(112, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(112, 18).Value2 = "Provincia de Cardenal Caro"
$ws.Cells.Item(112, 19).Value2 = 2500
$ws.Cells.Item(112, 20).Value2 = 2
